# Apply edits to the "s2" worksheet:
#  - Rename header cells A1/B1 from "asdasd"/"cadabra" to "column_1"/"column_2"
#  - Update B9 from "one" to "two"
#  - Update B10 from "two" to "three "
#  - Move the active selection to D16

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("s2")
$ws.Activate()

$ws.Range("A1").Value = "column_1"
$ws.Range("B1").Value = "column_2"

$ws.Range("B9").Value = "two"
$ws.Range("B10").Value = "three "

$ws.Range("D16").Select()
